$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Update the "Locations: {{org_name}}" label to "School: {{org_name}}" ---
# (dashboard filter is now by school instead of by location), keeping the
# existing rich-text formatting: the "School: " prefix stays bold, the
# "{{org_name}}" placeholder stays regular.
$cell1 = $ws1.Range("B4")
$cell1.Value = "School: {{org_name}}"
$boldPart = $cell1.Characters(1, 8)
$boldPart.Font.Bold = $true
$boldPart.Font.Name = "Times New Roman"
$boldPart.Font.Size = 12
$restPart = $cell1.Characters(9, 12)
$restPart.Font.Bold = $false
$restPart.Font.Name = "Times New Roman"
$restPart.Font.Size = 12

# Mirror the same label onto the other sheet; copy/paste so both sheets
# end up sharing a single (deduplicated) rich-text string entry.
$cell1.Copy() | Out-Null
$ws2.Range("B4").PasteSpecial(-4104) | Out-Null
$excel.CutCopyMode = $false

# --- Switch the active sheet / selection to the "Online" dashboard ---
$ws1.Range("E5").Select() | Out-Null
$ws2.Activate() | Out-Null
$ws2.Range("B4").Select() | Out-Null
